$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.203125
$ws.Range("C2").Value = 0.51875
$ws.Range("J2").Value = 0.009375
$ws.Range("P2").Value = 0.175
$ws.Range("S2").Value = 0.09375
# Row 3
$ws.Range("B3").Value = 0.01129943502824859
$ws.Range("C3").Value = 0.01694915254237288
$ws.Range("J3").Value = 0.02259887005649718
$ws.Range("P3").Value = 0.7457627118644068
$ws.Range("S3").Value = 0.2033898305084746
# Row 4
$ws.Range("J4").Value = 0.04444444444444445
$ws.Range("P4").Value = 0.5333333333333333
$ws.Range("S4").Value = 0.4222222222222222
# Row 6
$ws.Range("B6").Value = 0.06143344709897611
$ws.Range("D6").Value = 0.0136518771331058
$ws.Range("F6").Value = 0.04778156996587031
$ws.Range("J6").Value = 0.2866894197952218
$ws.Range("O6").Value = 0.02389078498293516
$ws.Range("Q6").Value = 0.1433447098976109
$ws.Range("R6").Value = 0.06143344709897611
$ws.Range("S6").Value = 0.3617747440273038
# Row 7
$ws.Range("B7").Value = 0.1697247706422018
$ws.Range("D7").Value = 0.009174311926605505
$ws.Range("F7").Value = 0.05963302752293578
$ws.Range("J7").Value = 0.1422018348623853
$ws.Range("O7").Value = 0.04587155963302753
$ws.Range("Q7").Value = 0.1192660550458716
$ws.Range("R7").Value = 0.04128440366972477
$ws.Range("S7").Value = 0.4128440366972477
# Row 8
$ws.Range("B8").Value = 0.09394572025052192
$ws.Range("D8").Value = 0.01878914405010438
$ws.Range("F8").Value = 0.05219206680584551
$ws.Range("J8").Value = 0.1315240083507307
$ws.Range("O8").Value = 0.04384133611691023
$ws.Range("Q8").Value = 0.1377870563674322
$ws.Range("R8").Value = 0.06263048016701461
$ws.Range("S8").Value = 0.4592901878914405
# Row 9
$ws.Range("B9").Value = 0.1038961038961039
$ws.Range("D9").Value = 0.02164502164502164
$ws.Range("E9").Value = 0.004329004329004329
$ws.Range("F9").Value = 0.09523809523809523
$ws.Range("J9").Value = 0.1125541125541126
$ws.Range("O9").Value = 0.03896103896103896
$ws.Range("Q9").Value = 0.1341991341991342
$ws.Range("R9").Value = 0.08658008658008658
$ws.Range("S9").Value = 0.4025974025974026
# Row 10
$ws.Range("B10").Value = 0.1226645004061738
$ws.Range("D10").Value = 0.02274573517465475
$ws.Range("E10").Value = 0.0008123476848090983
$ws.Range("F10").Value = 0.0901705930138099
$ws.Range("J10").Value = 0.1031681559707555
$ws.Range("O10").Value = 0.03330625507717303
$ws.Range("Q10").Value = 0.1933387489845654
$ws.Range("R10").Value = 0.0503655564581641
$ws.Range("S10").Value = 0.3834281072298944
# Row 11
$ws.Range("G11").Value = 0.1629392971246006
$ws.Range("J11").Value = 0.08626198083067092
$ws.Range("K11").Value = 0.1916932907348243
$ws.Range("L11").Value = 0.5271565495207667
$ws.Range("S11").Value = 0.03194888178913738
# Row 12
$ws.Range("G12").Value = 0.6832579185520362
$ws.Range("J12").Value = 0.1447963800904978
$ws.Range("K12").Value = 0.02714932126696833
$ws.Range("L12").Value = 0.06334841628959276
$ws.Range("S12").Value = 0.08144796380090498
# Row 13
$ws.Range("F13").Value = 0.0196078431372549
$ws.Range("G13").Value = 0.7058823529411765
$ws.Range("J13").Value = 0.09803921568627451
$ws.Range("S13").Value = 0.1764705882352941
# Row 15
$ws.Range("F15").Value = 0.04395604395604396
$ws.Range("H15").Value = 0.1721611721611722
$ws.Range("I15").Value = 0.06227106227106227
$ws.Range("J15").Value = 0.282051282051282
$ws.Range("K15").Value = 0.0695970695970696
$ws.Range("M15").Value = 0.007326007326007326
$ws.Range("N15").Value = 0.003663003663003663
$ws.Range("O15").Value = 0.05860805860805861
$ws.Range("S15").Value = 0.3003663003663004
# Row 16
$ws.Range("F16").Value = 0.03883495145631068
$ws.Range("H16").Value = 0.1990291262135922
$ws.Range("I16").Value = 0.09223300970873786
$ws.Range("J16").Value = 0.3446601941747573
$ws.Range("K16").Value = 0.1213592233009709
$ws.Range("M16").Value = 0.02427184466019417
$ws.Range("N16").Value = 0.004854368932038835
$ws.Range("O16").Value = 0.03883495145631068
$ws.Range("S16").Value = 0.1359223300970874
# Row 17
$ws.Range("F17").Value = 0.0225
$ws.Range("H17").Value = 0.155
$ws.Range("I17").Value = 0.095
$ws.Range("J17").Value = 0.4275
$ws.Range("K17").Value = 0.1
$ws.Range("M17").Value = 0.0225
$ws.Range("O17").Value = 0.0725
$ws.Range("S17").Value = 0.105
# Row 18
$ws.Range("F18").Value = 0.02142857142857143
$ws.Range("H18").Value = 0.1571428571428571
$ws.Range("I18").Value = 0.1142857142857143
$ws.Range("J18").Value = 0.3928571428571428
$ws.Range("K18").Value = 0.09285714285714286
$ws.Range("M18").Value = 0.01428571428571429
$ws.Range("O18").Value = 0.07142857142857142
$ws.Range("S18").Value = 0.1357142857142857
# Row 19
$ws.Range("F19").Value = 0.04134548002803083
$ws.Range("H19").Value = 0.2032235459004905
$ws.Range("I19").Value = 0.08969866853538892
$ws.Range("J19").Value = 0.3419761737911703
$ws.Range("K19").Value = 0.09600560616678346
$ws.Range("M19").Value = 0.01751927119831815
$ws.Range("N19").Value = 0.002102312543798178
$ws.Range("O19").Value = 0.07077785564120533
$ws.Range("S19").Value = 0.1373510861948143

Write-Output "Applied team specific matrix updates"